$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = newPrice; E = newVolume } for cells that changed,
# taken from the crypto-ticker refresh commit.
$changes = @{
    2  = @{ D = "67.840.57"; E = "  +0.79%  " }
    3  = @{ D = "2.624.73";  E = "  +0.22%  " }
    4  = @{ E = "  +0.02%  " }
    5  = @{ D = "596.24";    E = "  +0.03%  " }
    6  = @{ D = "153.53";    E = "  +0.35%  " }
    7  = @{ E = "  +0.00%  " }
    8  = @{ D = "0.549";     E = "  -1.43%  " }
    9  = @{ D = "2.623.77";  E = "  +0.24%  " }
    10 = @{ E = "  +9.38%  " }
    11 = @{ E = "  -0.84%  " }
    12 = @{ E = "  +0.58%  " }
    13 = @{ E = "  -0.21%  " }
    14 = @{ D = "27.54";     E = "  -0.75%  " }
    15 = @{ E = "  +3.83%  " }
    16 = @{ D = "3.104.63";  E = "  +0.30%  " }
    17 = @{ D = "67.864.09"; E = "  +1.05%  " }
    18 = @{ D = "2.619.60";  E = "  +0.11%  " }
    19 = @{ D = "11.39";     E = "  +2.64%  " }
    20 = @{ D = "370.72";    E = "  +1.99%  " }
    21 = @{ E = "  -0.10%  " }
    22 = @{ D = "4.25";      E = "  -1.07%  " }
    23 = @{ E = "  -2.14%  " }
    24 = @{ E = "  -1.36%  " }
    25 = @{ D = "72.04";     E = "  +1.34%  " }
    26 = @{ D = "1.00";      E = "  +0.01%  " }
    27 = @{ D = "9.92";      E = "  -1.04%  " }
    29 = @{ E = "  +2.00%  " }
    30 = @{ E = "  -0.03%  " }
    31 = @{ D = "576.95";    E = "  -1.57%  " }
    32 = @{ D = "7.88";      E = "  +0.64%  " }
    33 = @{ E = "  -0.69%  " }
    34 = @{ E = "  +0.70%  " }
    35 = @{ D = "1.00";      E = "  +0.03%  " }
    36 = @{ E = "  +0.32%  " }
    37 = @{ E = "  +0.20%  " }
    38 = @{ D = "159.20";    E = "  +1.15%  " }
    39 = @{ D = "19.14";     E = "  +0.30%  " }
    40 = @{ E = "  +4.55%  " }
    41 = @{ D = "0.367";     E = "  -0.19%  " }
    42 = @{ E = "  +1.24%  " }
    43 = @{ D = "0.0₆0332"; E = "  +14.86%  " }
    44 = @{ D = "2.63";      E = "  +2.38%  " }
    45 = @{ D = "17.39";     E = "  +6.23%  " }
    46 = @{ E = "  +0.06%  " }
    47 = @{ D = "40.18";     E = "  -2.33%  " }
    48 = @{ D = "155.54";    E = "  -0.78%  " }
    49 = @{ E = "  -1.17%  " }
    50 = @{ D = "21.84";     E = "  -0.60%  " }
    51 = @{ E = "  -1.22%  " }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]

    if ($cols.ContainsKey("D")) {
        $cell = $ws.Cells.Item($row, 4)
        # Preserve the original style while forcing a text number-format so
        # that numeric-looking strings (e.g. "596.24", "1.00") are not
        # silently coerced into real numbers by Excel's auto-detection.
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $cols["D"]
        $cell.Style = $origStyle
    }

    if ($cols.ContainsKey("E")) {
        $ws.Cells.Item($row, 5).Value = $cols["E"]
    }
}
